{"js": "// The document contains four transcription records, each with an\n// \"<id>...</id>\" marker split across three runs:\n//   run 1: \"<id>\"      (Courier New, color 7F6000, sz 18)\n//   run 2: \"p047v_aN\"  (default font, color 000000)\n//   run 3: \"</id>\"     (Courier New, color 7F6000, sz 18)\n//\n// The edit renumbers the ids (dropping the \"a\") and merges the three\n// runs into a single run carrying the \"<id>\" run's formatting, e.g.\n// \"<id>p047v_a1</id>\" -> \"<id>p047v_1</id>\" as one run.\n//\n// Search across the whole visible string (it spans run boundaries)\n// and replace it in one shot; Word collapses the hit into a single\n// run using the formatting of the first character of the match.\nconst body = context.document.body;\n\nfor (let n = 1; n <= 4; n++) {\n  const oldText = `<id>p047v_a${n}</id>`;\n  const newText = `<id>p047v_${n}</id>`;\n\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document holds four transcription records, each with an\n# \"<id>...</id>\" marker split across three runs:\n#   run 1: \"<id>\"      (Courier New, color 7F6000, sz 18)\n#   run 2: \"p047v_aN\"  (default font, color 000000)\n#   run 3: \"</id>\"     (Courier New, color 7F6000, sz 18)\n#\n# The edit renumbers the ids (dropping the \"a\") and collapses the\n# three runs into a single run carrying the \"<id>\" run's formatting,\n# e.g. \"<id>p047v_a1</id>\" -> \"<id>p047v_1</id>\" as one run.\n#\n# Find/Replace across the whole visible string (it spans the run\n# boundaries); Word merges the hit into a single run that takes on\n# the formatting of the first character of the match.\n$d = $word.ActiveDocument\n\nfor ($n = 1; $n -le 4; $n++) {\n    $oldText = \"<id>p047v_a$n</id>\"\n    $newText = \"<id>p047v_$n</id>\"\n\n    $rng = $d.Content\n    $rng.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
